$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F9: 208 -> 207; G9: 6150.56 -> 6120.99
$ws.Range("F9").Value = 207
$ws.Range("G9").Value = 6120.99

# B10: 46121.3 -> 46091.73
$ws.Range("B10").Value = 46091.73

# F34: 80 -> 79; G34: 2457.6 -> 2426.88
$ws.Range("F34").Value = 79
$ws.Range("G34").Value = 2426.88

# B47: 19482.89 -> 19452.17
$ws.Range("B47").Value = 19452.17

# F83: 82 -> 81; G83: 5476.78 -> 5409.99
$ws.Range("F83").Value = 81
$ws.Range("G83").Value = 5409.99

# F85: 25 -> 15; G85: 568.25 -> 340.95
$ws.Range("F85").Value = 15
$ws.Range("G85").Value = 340.95

# F102: 33 -> 30; G102: 3221.46 -> 2928.6
$ws.Range("F102").Value = 30
$ws.Range("G102").Value = 2928.6

# F104: 25 -> 23; G104: 2561.5 -> 2356.58
$ws.Range("F104").Value = 23
$ws.Range("G104").Value = 2356.58

# F106: 144 -> 143; G106: 21696.48 -> 21545.81
$ws.Range("F106").Value = 143
$ws.Range("G106").Value = 21545.81

# F109: 165 -> 164; G109: 20702.55 -> 20577.08
$ws.Range("F109").Value = 164
$ws.Range("G109").Value = 20577.08

# B114: 320411.49 -> 319343.48
$ws.Range("B114").Value = 319343.48

# B233: 48719 -> 64979; E233: 353.35 -> 314.41; F233: -81 -> 33; G233: -23955.75 -> 9759.75
$ws.Range("B233").Value = 64979
$ws.Range("E233").Value = 314.41
$ws.Range("F233").Value = 33
$ws.Range("G233").Value = 9759.75

# B234: 64979 -> 48719; E234: 314.41 -> 353.35; F234: 33 -> -81; G234: 9759.75 -> -23955.75
$ws.Range("B234").Value = 48719
$ws.Range("E234").Value = 353.35
$ws.Range("F234").Value = -81
$ws.Range("G234").Value = -23955.75

# B294: 57802 -> 63531; E294: 162.71 -> 152.53; F294: -79 -> 80; G294: -11334.92 -> 11478.4
$ws.Range("B294").Value = 63531
$ws.Range("E294").Value = 152.53
$ws.Range("F294").Value = 80
$ws.Range("G294").Value = 11478.4

# B296: 63531 -> 57802; E296: 152.53 -> 162.71; F296: 80 -> -79; G296: 11478.4 -> -11334.92
$ws.Range("B296").Value = 57802
$ws.Range("E296").Value = 162.71
$ws.Range("F296").Value = -79
$ws.Range("G296").Value = -11334.92

# B311: 61605 -> 63563; E311: 133.78 -> 119.04; F311: -13 -> 2; G311: -1455.48 -> 223.92
$ws.Range("B311").Value = 63563
$ws.Range("E311").Value = 119.04
$ws.Range("F311").Value = 2
$ws.Range("G311").Value = 223.92

# B312: 63563 -> 61605; E312: 119.04 -> 133.78; F312: 2 -> -13; G312: 223.92 -> -1455.48
$ws.Range("B312").Value = 61605
$ws.Range("E312").Value = 133.78
$ws.Range("F312").Value = -13
$ws.Range("G312").Value = -1455.48

# B315: 63560 -> 60325; E315: 134.87 -> 151.57; F315: 1 -> -102; G315: 126.86 -> -12939.72
$ws.Range("B315").Value = 60325
$ws.Range("E315").Value = 151.57
$ws.Range("F315").Value = -102
$ws.Range("G315").Value = -12939.72

# B316: 60325 -> 63560; E316: 151.57 -> 134.87; F316: -102 -> 1; G316: -12939.72 -> 126.86
$ws.Range("B316").Value = 63560
$ws.Range("E316").Value = 134.87
$ws.Range("F316").Value = 1
$ws.Range("G316").Value = 126.86

# F319: 11 -> 2; G319: 131.34 -> 23.88
$ws.Range("F319").Value = 2
$ws.Range("G319").Value = 23.88

# F328: 1680 -> 1667; G328: 35330.4 -> 35057.01
$ws.Range("F328").Value = 1667
$ws.Range("G328").Value = 35057.01

# F333: 803 -> 802; G333: 137577.99 -> 137406.66
$ws.Range("F333").Value = 802
$ws.Range("G333").Value = 137406.66

# B339: 416337.8 -> 415785.62
$ws.Range("B339").Value = 415785.62

# F355: 133 -> 131; G355: 4287.92 -> 4223.44
$ws.Range("F355").Value = 131
$ws.Range("G355").Value = 4223.44

# B361: 16059.29 -> 15994.81
$ws.Range("B361").Value = 15994.81

# F364: 28 -> 27; G364: 3355.24 -> 3235.41
$ws.Range("F364").Value = 27
$ws.Range("G364").Value = 3235.41

# F389: 6 -> 5; G389: 335.16 -> 279.3
$ws.Range("F389").Value = 5
$ws.Range("G389").Value = 279.3

# B395: 281782.18 -> 281606.49
$ws.Range("B395").Value = 281606.49

# B420: 47097 -> 58047; D420: 112.28 -> 105.54; E420: 134.16 -> 126.1; F420: 15 -> 43; G420: 1684.2 -> 4538.22
$ws.Range("B420").Value = 58047
$ws.Range("D420").Value = 105.54
$ws.Range("E420").Value = 126.1
$ws.Range("F420").Value = 43
$ws.Range("G420").Value = 4538.22

# B421: 58047 -> 47097; D421: 105.54 -> 112.28; E421: 126.1 -> 134.16; F421: 43 -> 15; G421: 4538.22 -> 1684.2
$ws.Range("B421").Value = 47097
$ws.Range("D421").Value = 112.28
$ws.Range("E421").Value = 134.16
$ws.Range("F421").Value = 15
$ws.Range("G421").Value = 1684.2

# F435: 34 -> 33; G435: 1229.78 -> 1193.61
$ws.Range("F435").Value = 33
$ws.Range("G435").Value = 1193.61

# F446: 171 -> 170; G446: 10501.11 -> 10439.7
$ws.Range("F446").Value = 170
$ws.Range("G446").Value = 10439.7

# B448: 49249.42 -> 49151.84
$ws.Range("B448").Value = 49151.84

# B465: 65069 -> 53757; E465: 14.3 -> 16.08; F465: 23 -> -159; G465: 309.35 -> -2138.55
$ws.Range("B465").Value = 53757
$ws.Range("E465").Value = 16.08
$ws.Range("F465").Value = -159
$ws.Range("G465").Value = -2138.55

# B466: 53757 -> 65069; E466: 16.08 -> 14.3; F466: -159 -> 23; G466: -2138.55 -> 309.35
$ws.Range("B466").Value = 65069
$ws.Range("E466").Value = 14.3
$ws.Range("F466").Value = 23
$ws.Range("G466").Value = 309.35

# B472: 64915 -> 45695; E472: 20.98 -> 23.58; F472: 0 -> -36; G472: 0 -> -710.28
$ws.Range("B472").Value = 45695
$ws.Range("E472").Value = 23.58
$ws.Range("F472").Value = -36
$ws.Range("G472").Value = -710.28

# B473: 45695 -> 64915; E473: 23.58 -> 20.98; F473: -36 -> 0; G473: -710.28 -> 0
$ws.Range("B473").Value = 64915
$ws.Range("E473").Value = 20.98
$ws.Range("F473").Value = 0
$ws.Range("G473").Value = 0

# B479: 45718 -> 64927; E479: 19.38 -> 17.26; F479: -294 -> 256; G479: -4768.68 -> 4152.32
$ws.Range("B479").Value = 64927
$ws.Range("E479").Value = 17.26
$ws.Range("F479").Value = 256
$ws.Range("G479").Value = 4152.32

# B480: 64927 -> 45718; E480: 17.26 -> 19.38; F480: 256 -> -294; G480: 4152.32 -> -4768.68
$ws.Range("B480").Value = 45718
$ws.Range("E480").Value = 19.38
$ws.Range("F480").Value = -294
$ws.Range("G480").Value = -4768.68

# F511: 113 -> 111; G511: 2740.25 -> 2691.75
$ws.Range("F511").Value = 111
$ws.Range("G511").Value = 2691.75

# B528: 25067.81 -> 25019.31
$ws.Range("B528").Value = 25019.31

# F584: 99 -> 98; G584: 3250.17 -> 3217.34
$ws.Range("F584").Value = 98
$ws.Range("G584").Value = 3217.34

# B585: 64833 -> 60025; E585: 34.9 -> 37.22; F585: 96 -> -98; G585: 3151.68 -> -3217.34
$ws.Range("B585").Value = 60025
$ws.Range("E585").Value = 37.22
$ws.Range("F585").Value = -98
$ws.Range("G585").Value = -3217.34

# B586: 60025 -> 64833; E586: 37.22 -> 34.9; F586: -98 -> 96; G586: -3217.34 -> 3151.68
$ws.Range("B586").Value = 64833
$ws.Range("E586").Value = 34.9
$ws.Range("F586").Value = 96
$ws.Range("G586").Value = 3151.68

# F587: 148 -> 147; G587: 4858.84 -> 4826.01
$ws.Range("F587").Value = 147
$ws.Range("G587").Value = 4826.01

# B596: 64830 -> 60022; E596: 34.9 -> 37.22; F596: 114 -> -113; G596: 3742.62 -> -3709.79
$ws.Range("B596").Value = 60022
$ws.Range("E596").Value = 37.22
$ws.Range("F596").Value = -113
$ws.Range("G596").Value = -3709.79

# B597: 60022 -> 64830; E597: 37.22 -> 34.9; F597: -113 -> 113; G597: -3709.79 -> 3709.79
$ws.Range("B597").Value = 64830
$ws.Range("E597").Value = 34.9
$ws.Range("F597").Value = 113
$ws.Range("G597").Value = 3709.79

# B598: 47850.54 -> 47752.05
$ws.Range("B598").Value = 47752.05

# F701: 250 -> 248; G701: 35782.5 -> 35496.24
$ws.Range("F701").Value = 248
$ws.Range("G701").Value = 35496.24

# F704: 38 -> 37; G704: 5057.8 -> 4924.7
$ws.Range("F704").Value = 37
$ws.Range("G704").Value = 4924.7

# F707: 193 -> 192; G707: 4191.96 -> 4170.24
$ws.Range("F707").Value = 192
$ws.Range("G707").Value = 4170.24

# F715: 506 -> 505; G715: 61079.26 -> 60958.55
$ws.Range("F715").Value = 505
$ws.Range("G715").Value = 60958.55

# B716: 245794.38 -> 245232.59
$ws.Range("B716").Value = 245232.59

# F720: 88 -> 87; G720: 14422.32 -> 14258.43
$ws.Range("F720").Value = 87
$ws.Range("G720").Value = 14258.43

# F722: 14 -> 12; G722: 1523.34 -> 1305.72
$ws.Range("F722").Value = 12
$ws.Range("G722").Value = 1305.72

# B743: 131525.59 -> 131144.08
$ws.Range("B743").Value = 131144.08

# F768: 3763 -> 3757; G768: 613782.9300000001 -> 612804.27
$ws.Range("F768").Value = 3757
$ws.Range("G768").Value = 612804.27

# F769: 98 -> 97; G769: 17264.66 -> 17088.49
$ws.Range("F769").Value = 97
$ws.Range("G769").Value = 17088.49

# F772: 10 -> 9; G772: 771.4 -> 694.26
$ws.Range("F772").Value = 9
$ws.Range("G772").Value = 694.26

# B775: 926160.51 -> 924928.54
$ws.Range("B775").Value = 924928.54

# F789: 70 -> 65; G789: 2725.1 -> 2530.45
$ws.Range("F789").Value = 65
$ws.Range("G789").Value = 2530.45

# B792: 99715.37 -> 99520.72
$ws.Range("B792").Value = 99520.72

# B793: 3881927.51 -> 3877392.37
$ws.Range("B793").Value = 3877392.37

# B794: 3881927.51 -> 3877392.37
$ws.Range("B794").Value = 3877392.37

